# Start implementing missings in minimal example:
# Insert two new rows right after the header row (new rows 2 and 3),
# pushing the existing data rows down by two, and fill them with the
# new "varMetrisch" missing-value code definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:3").Insert()

$ws.Cells.Item(2, 1).Value = "varMetrisch"
$ws.Cells.Item(2, 2).Value = -99
$ws.Cells.Item(2, 3).Value = "ja"
$ws.Cells.Item(2, 4).Value = "not reached"
$ws.Cells.Item(2, 5).Value = "nein"

$ws.Cells.Item(3, 1).Value = "varMetrisch"
$ws.Cells.Item(3, 2).Value = -98
$ws.Cells.Item(3, 3).Value = "ja"
$ws.Cells.Item(3, 4).Value = "omission"
$ws.Cells.Item(3, 5).Value = "nein"
